# Updated cryptos list on Sat Dec 23 21:38:46 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped
# values. Some Price cells are plain numeric-looking strings (e.g. "270.50")
# that must stay stored as text (matching the sheet's original inline-string
# cells), so for those we briefly force a Text number format before writing
# the value and then restore the cell's original (Normal) style so no
# visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.769.99'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.291.22'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '102.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '270.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.618'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.62%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.607'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.89'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0934'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.11'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.58'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('E15').Value = '  -2.21%  '
$ws.Range('D16').Value = '2.287.47'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').Value = '43.758.76'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '233.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.76%  '
$ws.Range('E23').Value = '  +13.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '40.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.97%  '
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '177.62'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '21.82'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0902'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.69%  '
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.89'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.49%  '
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.54'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.14%  '
$ws.Range('E39').Value = '  -2.66%  '
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.20'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.29'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('E45').Value = '  -4.60%  '
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '99.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('E49').Value = '  +11.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.440'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.92%  '
$ws.Range('D51').Value = '2.523.36'
$ws.Range('E51').Value = '  -0.85%  '
